$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1243
$ws1.Range("F4").Value = 17069
$ws1.Range("G4").Value = 65
$ws1.Range("F5").Value = 45
$ws1.Range("F9").Value = 34
$ws1.Range("F12").Value = 133
$ws1.Range("F13").Value = 11803
$ws1.Range("F14").Value = 31
$ws1.Range("F15").Value = 33
$ws1.Range("F16").Value = 1484
$ws1.Range("F17").Value = 4688
$ws1.Range("F22").Value = 912
$ws1.Range("F23").Value = 345
$ws1.Range("F24").Value = 154
$ws1.Range("F25").Value = 40

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1243
$ws4.Range("F4").Value = 17069
$ws4.Range("G4").Value = 65
$ws4.Range("F5").Value = 45
$ws4.Range("F9").Value = 35
$ws4.Range("F12").Value = 133
$ws4.Range("F15").Value = 11803
$ws4.Range("F16").Value = 31
$ws4.Range("F17").Value = 33
$ws4.Range("F18").Value = 1484
$ws4.Range("F19").Value = 4688
$ws4.Range("F24").Value = 912
$ws4.Range("F25").Value = 345
$ws4.Range("F26").Value = 154
$ws4.Range("F27").Value = 40
